$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Doctor D001's appointment on 2024-10-30 (row 2): Available End Time 17:00 -> 10:00
$ws.Range("D2").Value = 0.41666666666666669

# Doctor D001's appointment on 2024-10-31 (row 3): Available Start Time 10:00 -> 14:00
$ws.Range("C3").Value = 0.58333333333333337

# Reflect where the user's selection ended up after creating the new appointment
$ws.Range("D10").Select() | Out-Null
